$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold text-like values (e.g. "219.70", "0.531"),
# so force text storage on any cell we actually write -- otherwise Excel
# would silently coerce numeric-looking strings (losing trailing zeros /
# dotted thousand separators) into real numbers.
function Set-TextValue($cell, $value) {
    if ($null -eq $value) { return }
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Simple per-row price / volume updates (rows where coin identity doesn't change)
$updates = @(
    @{ Row = 2;  D = "30.525.61";  E = "  +1.92%  " }
    @{ Row = 3;  D = "1.675.59";   E = "  +2.55%  " }
    @{ Row = 4;  D = $null;        E = "  +0.16%  " }
    @{ Row = 5;  D = "219.70";     E = "  +2.39%  " }
    @{ Row = 6;  D = "0.529";      E = "  +1.99%  " }
    @{ Row = 7;  D = $null;        E = "  +0.20%  " }
    @{ Row = 8;  D = "29.87";      E = "  +3.96%  " }
    @{ Row = 9;  D = $null;        E = "  +2.44%  " }
    @{ Row = 10; D = "0.0635";     E = "  +4.26%  " }
    @{ Row = 11; D = "0.0906";     E = "  -0.77%  " }
    @{ Row = 12; D = "1.919.21";   E = "  +2.76%  " }

    @{ Row = 16; D = "3.97";       E = "  +2.80%  " }
    @{ Row = 17; D = "30.556.44";  E = "  +2.00%  " }
    @{ Row = 18; D = "66.30";      E = "  +3.37%  " }
    @{ Row = 19; D = "244.37";     E = "  +0.47%  " }
    @{ Row = 20; D = "0.0₃0721";   E = "  +2.79%  " }
    @{ Row = 21; D = $null;        E = "  +0.04%  " }
    @{ Row = 22; D = "4.26";       E = "  +3.06%  " }
    @{ Row = 23; D = "10.01";      E = "  +1.32%  " }
    @{ Row = 24; D = $null;        E = "  +1.09%  " }
    @{ Row = 25; D = "158.02";     E = $null }
    @{ Row = 26; D = "15.91";      E = "  +2.38%  " }
    @{ Row = 27; D = "0.112";      E = "  +2.12%  " }
    @{ Row = 28; D = "6.68";       E = "  +1.09%  " }
    @{ Row = 29; D = $null;        E = "  +0.13%  " }
    @{ Row = 30; D = "0.0496";     E = "  +2.15%  " }
    @{ Row = 31; D = $null;        E = "  +2.99%  " }
    @{ Row = 32; D = "3.48";       E = "  +2.93%  " }
    @{ Row = 33; D = "1.508.57";   E = "  +5.92%  " }
    @{ Row = 34; D = "3.29";       E = "  +3.87%  " }
    @{ Row = 35; D = "1.76";       E = "  +7.24%  " }

    @{ Row = 38; D = "0.601";      E = "  +8.96%  " }
    @{ Row = 39; D = $null;        E = "  +5.44%  " }
    @{ Row = 40; D = "2.70";       E = "  -3.69%  " }
    @{ Row = 41; D = "2.31";       E = "  +0.54%  " }
    @{ Row = 42; D = "0.839";      E = "  +1.03%  " }
    @{ Row = 43; D = "1.98";       E = "  -0.58%  " }
    @{ Row = 44; D = "0.0498";     E = "  +1.80%  " }
    @{ Row = 45; D = $null;        E = "  +0.91%  " }
    @{ Row = 46; D = $null;        E = "  +0.13%  " }
    @{ Row = 47; D = "5.58";       E = "  +4.26%  " }
    @{ Row = 48; D = "51.24";      E = "  -3.27%  " }
    @{ Row = 49; D = "1.813.05";   E = $null }
    @{ Row = 50; D = "94.80";      E = "  +5.91%  " }
    @{ Row = 51; D = "0.0₆0112";   E = "  -0.39%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    Set-TextValue $ws.Cells.Item($r, 4) $u.D
    Set-TextValue $ws.Cells.Item($r, 5) $u.E
}

# Rows 13-15 were reshuffled (coin name/link rotated) with fresh price/volume figures.
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Cells.Item(13, 4) "1.690.66"
Set-TextValue $ws.Cells.Item(13, 5) "  +3.48%  "

$ws.Cells.Item(14, 2).Value = "Polygon"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Cells.Item(14, 4) "0.615"
Set-TextValue $ws.Cells.Item(14, 5) "  +9.14%  "

$ws.Cells.Item(15, 2).Value = "Chainlink"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Cells.Item(15, 4) "10.23"
Set-TextValue $ws.Cells.Item(15, 5) "  +10.05%  "

# Rows 36-37 were swapped (Aave <-> TrustWalletToken) with fresh price/volume figures.
$ws.Cells.Item(36, 2).Value = "TrustWalletToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Cells.Item(36, 4) "1.03"
Set-TextValue $ws.Cells.Item(36, 5) "  -0.38%  "

$ws.Cells.Item(37, 2).Value = "Aave"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Cells.Item(37, 4) "83.61"
Set-TextValue $ws.Cells.Item(37, 5) "  +10.19%  "
